{"js": "// Update the worksheet date title and the 25 division problems.\n// Values are addressed positionally (row/col in the single table, and the\n// first body paragraph for the title) rather than via text search, because\n// several problem cells share identical text (e.g. \"54\u00f79=\" appears twice)\n// and a global find/replace would not be able to distinguish them.\n\nconst body = context.document.body;\n\n// --- Title paragraph: \"2024-08-12 Monday\" -> \"2024-08-13 Tuesday\" ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\n\nif (titlePara.text === \"2024-08-12 Monday\") {\n  titlePara.getRange().insertText(\"2024-08-13 Tuesday\", \"Replace\");\n}\n\n// --- Table of division problems ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Old value -> new value, keyed by (row, col) in reading order of the\n// non-empty rows (0, 4, 8, 12, 16), each holding 5 cells.\nconst replacements = [\n  { row: 0, col: 0, oldText: \"31\u00f75=\", newText: \"53\u00f75=\" },\n  { row: 0, col: 1, oldText: \"85\u00f76=\", newText: \"79\u00f75=\" },\n  { row: 0, col: 2, oldText: \"69\u00f73=\", newText: \"43\u00f76=\" },\n  { row: 0, col: 3, oldText: \"44\u00f74=\", newText: \"79\u00f75=\" },\n  { row: 0, col: 4, oldText: \"73\u00f73=\", newText: \"84\u00f78=\" },\n\n  { row: 4, col: 0, oldText: \"91\u00f76=\", newText: \"29\u00f77=\" },\n  { row: 4, col: 1, oldText: \"52\u00f74=\", newText: \"13\u00f74=\" },\n  { row: 4, col: 2, oldText: \"14\u00f74=\", newText: \"88\u00f77=\" },\n  { row: 4, col: 3, oldText: \"19\u00f73=\", newText: \"55\u00f74=\" },\n  { row: 4, col: 4, oldText: \"74\u00f73=\", newText: \"64\u00f77=\" },\n\n  { row: 8, col: 0, oldText: \"35\u00f75=\", newText: \"21\u00f75=\" },\n  { row: 8, col: 1, oldText: \"67\u00f78=\", newText: \"64\u00f73=\" },\n  { row: 8, col: 2, oldText: \"45\u00f75=\", newText: \"45\u00f79=\" },\n  { row: 8, col: 3, oldText: \"32\u00f73=\", newText: \"43\u00f79=\" },\n  { row: 8, col: 4, oldText: \"23\u00f79=\", newText: \"76\u00f76=\" },\n\n  { row: 12, col: 0, oldText: \"54\u00f74=\", newText: \"43\u00f79=\" },\n  { row: 12, col: 1, oldText: \"29\u00f72=\", newText: \"30\u00f76=\" },\n  { row: 12, col: 2, oldText: \"23\u00f76=\", newText: \"81\u00f79=\" },\n  { row: 12, col: 3, oldText: \"80\u00f77=\", newText: \"77\u00f76=\" },\n  { row: 12, col: 4, oldText: \"57\u00f74=\", newText: \"56\u00f77=\" },\n\n  { row: 16, col: 0, oldText: \"79\u00f79=\", newText: \"85\u00f73=\" },\n  { row: 16, col: 1, oldText: \"76\u00f76=\", newText: \"57\u00f74=\" },\n  { row: 16, col: 2, oldText: \"54\u00f79=\", newText: \"69\u00f77=\" },\n  { row: 16, col: 3, oldText: \"12\u00f78=\", newText: \"58\u00f74=\" },\n  { row: 16, col: 4, oldText: \"54\u00f79=\", newText: \"95\u00f74=\" },\n];\n\nconst cells = replacements.map((r) => table.getCell(r.row, r.col));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const cell = cells[i];\n  const expected = replacements[i].oldText;\n  if (cell.value === expected) {\n    cell.getRange().insertText(replacements[i].newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date title and the 25 division problems.\n#\n# Values are addressed positionally (row/col of the single table, and the\n# document's first paragraph for the title) rather than via Find/Replace,\n# because several problem cells share identical text (e.g. \"54\u00f79=\" appears\n# twice in the original) and a global replace could not tell them apart.\n\n$d = $word.ActiveDocument\n\n# --- Title paragraph: \"2024-08-12 Monday\" -> \"2024-08-13 Tuesday\" ---\n$titlePara = $d.Paragraphs.Item(1)\n$titleRange = $titlePara.Range\n$titleText = $titleRange.Text -replace \"[\\r\\a]\", \"\"\nif ($titleText -eq \"2024-08-12 Monday\") {\n    $titleRange.Text = \"2024-08-13 Tuesday\"\n}\n\n# --- Table of division problems ---\n$table = $d.Tables.Item(1)\n\n# Each entry: 1-based table row/column, expected old text, new text.\n# Non-empty rows are 1, 5, 9, 13, 17 (each holding 5 problem cells);\n# the rows in between are blank spacer rows.\n$replacements = @(\n    @{ Row = 1;  Col = 1; Old = \"31\u00f75=\"; New = \"53\u00f75=\" },\n    @{ Row = 1;  Col = 2; Old = \"85\u00f76=\"; New = \"79\u00f75=\" },\n    @{ Row = 1;  Col = 3; Old = \"69\u00f73=\"; New = \"43\u00f76=\" },\n    @{ Row = 1;  Col = 4; Old = \"44\u00f74=\"; New = \"79\u00f75=\" },\n    @{ Row = 1;  Col = 5; Old = \"73\u00f73=\"; New = \"84\u00f78=\" },\n\n    @{ Row = 5;  Col = 1; Old = \"91\u00f76=\"; New = \"29\u00f77=\" },\n    @{ Row = 5;  Col = 2; Old = \"52\u00f74=\"; New = \"13\u00f74=\" },\n    @{ Row = 5;  Col = 3; Old = \"14\u00f74=\"; New = \"88\u00f77=\" },\n    @{ Row = 5;  Col = 4; Old = \"19\u00f73=\"; New = \"55\u00f74=\" },\n    @{ Row = 5;  Col = 5; Old = \"74\u00f73=\"; New = \"64\u00f77=\" },\n\n    @{ Row = 9;  Col = 1; Old = \"35\u00f75=\"; New = \"21\u00f75=\" },\n    @{ Row = 9;  Col = 2; Old = \"67\u00f78=\"; New = \"64\u00f73=\" },\n    @{ Row = 9;  Col = 3; Old = \"45\u00f75=\"; New = \"45\u00f79=\" },\n    @{ Row = 9;  Col = 4; Old = \"32\u00f73=\"; New = \"43\u00f79=\" },\n    @{ Row = 9;  Col = 5; Old = \"23\u00f79=\"; New = \"76\u00f76=\" },\n\n    @{ Row = 13; Col = 1; Old = \"54\u00f74=\"; New = \"43\u00f79=\" },\n    @{ Row = 13; Col = 2; Old = \"29\u00f72=\"; New = \"30\u00f76=\" },\n    @{ Row = 13; Col = 3; Old = \"23\u00f76=\"; New = \"81\u00f79=\" },\n    @{ Row = 13; Col = 4; Old = \"80\u00f77=\"; New = \"77\u00f76=\" },\n    @{ Row = 13; Col = 5; Old = \"57\u00f74=\"; New = \"56\u00f77=\" },\n\n    @{ Row = 17; Col = 1; Old = \"79\u00f79=\"; New = \"85\u00f73=\" },\n    @{ Row = 17; Col = 2; Old = \"76\u00f76=\"; New = \"57\u00f74=\" },\n    @{ Row = 17; Col = 3; Old = \"54\u00f79=\"; New = \"69\u00f77=\" },\n    @{ Row = 17; Col = 4; Old = \"12\u00f78=\"; New = \"58\u00f74=\" },\n    @{ Row = 17; Col = 5; Old = \"54\u00f79=\"; New = \"95\u00f74=\" }\n)\n\nforeach ($item in $replacements) {\n    $cell = $table.Cell($item.Row, $item.Col)\n    $cellRange = $cell.Range\n    $cellText = $cellRange.Text -replace \"[\\r\\a]\", \"\"\n    if ($cellText -eq $item.Old) {\n        $cellRange.Text = $item.New\n    }\n}\n"}
